$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.21   # Current Capital
$summary.Range("B4").Value = 0.2       # Total P&L $
$summary.Range("B5").Value = 0.06      # Total P&L %
$summary.Range("B6").Value = 70        # Total Trades
$summary.Range("B7").Value = 22        # Winning Trades
$summary.Range("B9").Value = 31.43     # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.21     # Capital
$status.Range("D4").Value = 70         # Trades
$status.Range("E4").Value = 0.2        # P&L $
$status.Range("F4").Value = 0.21       # P&L %
$status.Range("G4").Value = 31.43      # Win Rate %

# --- New trade row (#70, displayed as 70 in column A) appended to both
#     "All Trades" and "MarketMaking" sheets ---
$newRow = 71

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 70

    # Keep the date/time columns as plain text (matching the existing rows)
    # instead of letting Excel auto-convert them to date/time serials.
    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = "2026-02-17"

    $ws.Cells.Item($newRow, 3).NumberFormat = "@"
    $ws.Cells.Item($newRow, 3).Value = "15:48:15"

    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"
    $ws.Cells.Item($newRow, 5).Value = "UP"
    $ws.Cells.Item($newRow, 6).Value = 0.67
    $ws.Cells.Item($newRow, 7).Value = 0.72
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"
    $ws.Cells.Item($newRow, 9).Value = 7.4627
    $ws.Cells.Item($newRow, 10).Value = 0.05
    $ws.Cells.Item($newRow, 11).Value = 100.21
    $ws.Cells.Item($newRow, 12).Value = 0
    $ws.Cells.Item($newRow, 13).Value = 0
    $ws.Cells.Item($newRow, 14).Value = 0.6
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($newRow, 16).Value = "early_exit"
    $ws.Cells.Item($newRow, 17).Value = 0.15
}
